# Auto-generated Excel COM-interop script
# Applies the updated cryptocurrency "Price" (column D) and "Volume(1h)"
# (column E) figures to the active worksheet, matching the target diff.
#
# Several Price values are plain decimal-looking strings (e.g. "305.49").
# A naive ".Value = " assignment lets Excel's type-inference coerce those
# into floating point numbers (losing the original text formatting and
# introducing binary rounding noise, e.g. 305.49000000000001). To keep
# them as literal text -- matching the workbook's inlineStr/shared-string
# cells -- each such cell is briefly switched to a text number format
# before the value is written, then restored to the default "Normal"
# style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '44.023.78'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '2.237.83'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('E4').Value = '  +0.09%  '
Set-TextValue 'D5' '305.49'
$ws.Range('E5').Value = '  -3.97%  '
Set-TextValue 'D6' '95.31'
$ws.Range('E6').Value = '  -5.60%  '
Set-TextValue 'D7' '0.570'
$ws.Range('E7').Value = '  -1.48%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('E9').Value = '  -5.33%  '
Set-TextValue 'D10' '34.98'
$ws.Range('E10').Value = '  -5.89%  '
$ws.Range('E11').Value = '  -3.37%  '
Set-TextValue 'D12' '7.21'
$ws.Range('E12').Value = '  -5.01%  '
$ws.Range('E13').Value = '  -2.36%  '
$ws.Range('D14').Value = '2.576.41'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('D15').Value = '2.237.66'
$ws.Range('E15').Value = '  -0.99%  '
Set-TextValue 'D16' '0.825'
$ws.Range('E16').Value = '  -3.84%  '
$ws.Range('E17').Value = '  -6.41%  '
$ws.Range('D18').Value = '43.881.91'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('E19').Value = '  -2.55%  '
Set-TextValue 'D20' '12.22'
$ws.Range('E20').Value = '  -8.58%  '
Set-TextValue 'D21' '6.22'
$ws.Range('E21').Value = '  -3.76%  '
Set-TextValue 'D22' '64.94'
$ws.Range('E22').Value = '  -1.09%  '
Set-TextValue 'D23' '236.34'
$ws.Range('E23').Value = '  +0.63%  '
Set-TextValue 'D24' '2.93'
$ws.Range('E24').Value = '  -5.90%  '
Set-TextValue 'D25' '1.95'
$ws.Range('E25').Value = '  -5.78%  '
$ws.Range('E26').Value = '  +0.09%  '
Set-TextValue 'D27' '9.97'
$ws.Range('E27').Value = '  -6.36%  '
Set-TextValue 'D28' '37.94'
$ws.Range('E28').Value = '  -1.90%  '
$ws.Range('E29').Value = '  -1.68%  '
Set-TextValue 'D30' '5.96'
$ws.Range('E30').Value = '  -3.53%  '
Set-TextValue 'D31' '19.89'
$ws.Range('E31').Value = '  -1.48%  '
Set-TextValue 'D32' '154.76'
$ws.Range('E32').Value = '  -4.14%  '
$ws.Range('E33').Value = '  -5.31%  '
$ws.Range('E34').Value = '  +3.67%  '
$ws.Range('E35').Value = '  -3.37%  '
Set-TextValue 'D36' '0.119'
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('E37').Value = '  -5.55%  '
$ws.Range('E38').Value = '  -10.75%  '
Set-TextValue 'D39' '15.21'
$ws.Range('E39').Value = '  -8.35%  '
Set-TextValue 'D40' '3.36'
$ws.Range('E40').Value = '  -8.06%  '
$ws.Range('E41').Value = '  -8.67%  '
$ws.Range('E42').Value = '  -4.67%  '
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').Value = '1.739.92'
$ws.Range('E44').Value = '  -2.13%  '
Set-TextValue 'D45' '85.70'
$ws.Range('E45').Value = '  +5.26%  '
$ws.Range('E46').Value = '  -4.60%  '
Set-TextValue 'D47' '99.88'
$ws.Range('E47').Value = '  -4.25%  '
$ws.Range('E48').Value = '  -5.68%  '
Set-TextValue 'D49' '69.03'
$ws.Range('E49').Value = '  -7.43%  '
Set-TextValue 'D50' '8.08'
$ws.Range('E50').Value = '  -2.72%  '
Set-TextValue 'D51' '54.24'
$ws.Range('E51').Value = '  -6.52%  '
